# Add data for 2024-03-06: updates the 2024 (column K) crime-count figures
# across the citywide-totals, by-neighborhood summary, and the individual
# neighborhood sheets to reflect the newly added day's records (plus a
# handful of 2022/2023 reclassification corrections in columns I/J).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1256
$ws.Range("K3").Value = 1173
$ws.Range("I4").Value = 1784
$ws.Range("J4").Value = 1795
$ws.Range("K4").Value = 254
$ws.Range("K5").Value = 72
$ws.Range("J6").Value = 11060
$ws.Range("K6").Value = 1523
$ws.Range("I7").Value = 26237
$ws.Range("J7").Value = 29263
$ws.Range("K7").Value = 4278

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 250

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 41
$ws.Range("K3").Value = 25
$ws.Range("J4").Value = 23
$ws.Range("K6").Value = 13
$ws.Range("J7").Value = 591
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 51
$ws.Range("K3").Value = 67
$ws.Range("K4").Value = 10
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J6").Value = 262
$ws.Range("K6").Value = 53
$ws.Range("J7").Value = 902
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 23
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K8").Value = 250
$ws.Range("K11").Value = 92
$ws.Range("K14").Value = 25
$ws.Range("K15").Value = 35
$ws.Range("K16").Value = 9
$ws.Range("K18").Value = 33
$ws.Range("K19").Value = 114
$ws.Range("K25").Value = 20
$ws.Range("K29").Value = 197
$ws.Range("K31").Value = 48
$ws.Range("K33").Value = 177
$ws.Range("K36").Value = 47
$ws.Range("J37").Value = 902
$ws.Range("K37").Value = 138
$ws.Range("K42").Value = 139
$ws.Range("K43").Value = 41
$ws.Range("K46").Value = 9
$ws.Range("J48").Value = 324
$ws.Range("K48").Value = 44
$ws.Range("K49").Value = 29
$ws.Range("K52").Value = 119
$ws.Range("K53").Value = 57
$ws.Range("K55").Value = 45
$ws.Range("I63").Value = 193
$ws.Range("J63").Value = 89
$ws.Range("K63").Value = 15
$ws.Range("K64").Value = 30
$ws.Range("K65").Value = 110
$ws.Range("K66").Value = 18
$ws.Range("K67").Value = 179
$ws.Range("K72").Value = 17
$ws.Range("K76").Value = 57
$ws.Range("J83").Value = 591
$ws.Range("K83").Value = 82
$ws.Range("K85").Value = 213
$ws.Range("K94").Value = 53
$ws.Range("K95").Value = 76
$ws.Range("K96").Value = 63
$ws.Range("K97").Value = 33
$ws.Range("K99").Value = 79
$ws.Range("I101").Value = 26237
$ws.Range("J101").Value = 29263
$ws.Range("K101").Value = 4278

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 22
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 53
$ws.Range("K3").Value = 54
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 6
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 59
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J3").Value = 63
$ws.Range("K3").Value = 8
$ws.Range("K6").Value = 19
$ws.Range("J7").Value = 324
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 12
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 12
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 13
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 9
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 20
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 81
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 213

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 29
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 9
